{"js": "// Apply the pt_base v1.1.0 text revisions to the summary table.\n// Each edit rewrites only the text run content of a specific table\n// cell's first paragraph, preserving all existing run/paragraph\n// formatting by replacing the paragraph's range text in place.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nasync function setCellText(rowIndex, colIndex, newText) {\n  const cell = table.getCell(rowIndex, colIndex);\n  const paragraphs = cell.body.paragraphs;\n  paragraphs.load(\"items\");\n  await context.sync();\n  const range = paragraphs.items[0].getRange();\n  range.insertText(newText, Word.InsertLocation.replace);\n}\n\n// Header row: \"Summary 2\" -> \"Summary 0\"\nawait setCellText(0, 3, \"Summary 0\");\n\n// First results block (rows 1-4, column 0 labels)\nawait setCellText(1, 0, \"Age - mean (sd) [N]  \");\nawait setCellText(2, 0, \"Quality of life - median (IQR) [N]  \");\nawait setCellText(3, 0, \"Female - n (%) [N]  \");\nawait setCellText(4, 0, \"Ethnicity - n (%) [N]  \");\n\n// Second results block (rows 11-14)\nawait setCellText(11, 0, \"Age - mean (sd)  \");\nawait setCellText(11, 2, \"44.9 (10.1) \");\nawait setCellText(11, 3, \"44.6 (10.1) \");\nawait setCellText(11, 4, \"44.8 (10.1) \");\n\nawait setCellText(12, 0, \"Quality of life - median (IQR) [N (%)]  \");\n\nawait setCellText(13, 0, \"Female - n (%)  \");\nawait setCellText(13, 2, \"258 (51.0) \");\nawait setCellText(13, 3, \"261 (52.8) \");\nawait setCellText(13, 4, \"519 (51.9) \");\n\nawait setCellText(14, 0, \"Ethnicity - n (%) [N (%)]  \");\n\nawait context.sync();\n", "ps1": "# Apply the pt_base v1.1.0 text revisions to the summary table.\n# Each edit rewrites only the visible text of a specific table cell\n# (Cell.Range.Text), which keeps the cell's existing run/paragraph\n# formatting intact and only replaces the textual content.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Header row: \"Summary 2\" -> \"Summary 0\"\n$t.Cell(1, 4).Range.Text = \"Summary 0\"\n\n# First results block (rows 2-5, column 1 labels)\n$t.Cell(2, 1).Range.Text = \"Age - mean (sd) [N]  \"\n$t.Cell(3, 1).Range.Text = \"Quality of life - median (IQR) [N]  \"\n$t.Cell(4, 1).Range.Text = \"Female - n (%) [N]  \"\n$t.Cell(5, 1).Range.Text = \"Ethnicity - n (%) [N]  \"\n\n# Second results block (rows 12-15)\n$t.Cell(12, 1).Range.Text = \"Age - mean (sd)  \"\n$t.Cell(12, 3).Range.Text = \"44.9 (10.1) \"\n$t.Cell(12, 4).Range.Text = \"44.6 (10.1) \"\n$t.Cell(12, 5).Range.Text = \"44.8 (10.1) \"\n\n$t.Cell(13, 1).Range.Text = \"Quality of life - median (IQR) [N (%)]  \"\n\n$t.Cell(14, 1).Range.Text = \"Female - n (%)  \"\n$t.Cell(14, 3).Range.Text = \"258 (51.0) \"\n$t.Cell(14, 4).Range.Text = \"261 (52.8) \"\n$t.Cell(14, 5).Range.Text = \"519 (51.9) \"\n\n$t.Cell(15, 1).Range.Text = \"Ethnicity - n (%) [N (%)]  \"\n"}
